# Weekly update: a new "Ají" price record for Macroferia Regional de Talca was
# published. It belongs at the top of the existing data block (row 55, right
# after the header row and the most-recent-so-far row 54), pushing every
# existing record down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 55..143 down to 56..144, leaving row 55 free for the new record.
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A55").Value = 5
$ws.Range("B55").Value = "Macroferia Regional de Talca"
$ws.Range("C55").Value = "Maule"
$ws.Range("D55").Value = 44540
$ws.Range("E55").Value = 7
$ws.Range("F55").Value = 100112021
$ws.Range("G55").Value = "Ají"
$ws.Range("H55").Value = "Americana (o)"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 150
$ws.Range("K55").Value = 15000
$ws.Range("L55").Value = 15000
$ws.Range("M55").Value = 15000
$ws.Range("N55").Value = "`$/caja 15 kilos"
$ws.Range("O55").Value = "Región del Maule"
$ws.Range("P55").Value = 1000
$ws.Range("Q55").Value = 15
$ws.Range("R55").Value = "Hortaliza"
